$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (case with 380 kV), columns B-F, I-N
$data = @{
    2 = @{ 2 = 1.02; 3 = 1.04179778676742; 4 = 1.055958041975428; 5 = 1.050527668635203; 6 = 1.062669612866417; 9 = 1.042179269703827; 10 = 1.046877056667101; 11 = 1.058697071010305; 12 = 1.053281697198533; 13 = 1.065390332822686; 14 = 1.019419468148993 }
    3 = @{ 2 = 1.02; 3 = 1.043301607060659; 4 = 1.056857272790944; 5 = 1.051809633695733; 6 = 1.063913219580691; 9 = 1.042430696081661; 10 = 1.048024295307851; 11 = 1.059409627607577; 12 = 1.054374917207294; 13 = 1.06644772329165; 14 = 1.019807082982067 }
    4 = @{ 2 = 1.02; 3 = 1.044273375442894; 4 = 1.057436272003739; 5 = 1.052637882456593; 6 = 1.064715675237722; 9 = 1.042590352519559; 10 = 1.048764954894508; 11 = 1.059867208501922; 12 = 1.05508048722249; 13 = 1.067129126062385; 14 = 1.020057130080317 }
    5 = @{ 2 = 1.02; 3 = 1.044681601671565; 4 = 1.057678999674165; 5 = 1.05298577932928; 6 = 1.065052495149752; 9 = 1.042656746527012; 10 = 1.049075930769813; 11 = 1.060058741550729; 12 = 1.05537667831089; 13 = 1.067414921192908; 14 = 1.020162068211609 }
    6 = @{ 2 = 1.02; 3 = 1.044750126955782; 4 = 1.057719714684596; 5 = 1.053044175353784; 6 = 1.065109017588991; 9 = 1.042667851842725; 10 = 1.049128121824703; 11 = 1.060090851938607; 12 = 1.055426384974842; 13 = 1.067462868467641; 14 = 1.020179677179295 }
    7 = @{ 2 = 1.02; 3 = 1.044278831372414; 4 = 1.057439518026673; 5 = 1.052642532240344; 6 = 1.064720177929124; 9 = 1.042591242529839; 10 = 1.048769111724945; 11 = 1.059869771050766; 12 = 1.055084446628814; 13 = 1.067132947485818; 14 = 1.020058532980611 }
    8 = @{ 2 = 1.02; 3 = 1.042306284062583; 4 = 1.056262535241292; 5 = 1.050961180219622; 6 = 1.063090362067412; 9 = 1.04226486885355; 10 = 1.047265122268919; 11 = 1.058938606736776; 12 = 1.053651534733051; 13 = 1.065748264754625; 14 = 1.019550623903352 }
    9 = @{ 2 = 1.02; 3 = 1.03882010371349; 4 = 1.054166524582867; 5 = 1.047988494446116; 6 = 1.06020106921922; 9 = 1.041666496425991; 10 = 1.044601808093295; 11 = 1.057270957800778; 12 = 1.051112453850222; 13 = 1.063286669749315; 14 = 1.01864968622846 }
    10 = @{ 2 = 1.02; 3 = 1.036488616429037; 4 = 1.052754251629012; 5 = 1.045999723952488; 6 = 1.058262943910182; 9 = 1.041251897202484; 10 = 1.04281714878332; 11 = 1.056141038189564; 12 = 1.049409971064497; 13 = 1.061630840933837; 14 = 1.018044965684159 }
    11 = @{ 2 = 1.02; 3 = 1.035477209999127; 4 = 1.052139148972056; 5 = 1.045136839768178; 6 = 1.057420823631046; 9 = 1.041068639660371; 10 = 1.0420421383635; 11 = 1.055647435385866; 12 = 1.048670398774584; 13 = 1.060910291773757; 14 = 1.017782120167218 }
    12 = @{ 2 = 1.02; 3 = 1.035101240953575; 4 = 1.051910131833074; 5 = 1.044816059653054; 6 = 1.057107581646069; 9 = 1.041000007630824; 10 = 1.041753922049831; 11 = 1.055463434525219; 12 = 1.048395324667144; 13 = 1.060642107314674; 14 = 1.017684335705338 }
    13 = @{ 2 = 1.02; 3 = 1.035181900747484; 4 = 1.051959281298516; 5 = 1.044884880235113; 6 = 1.057174793160382; 9 = 1.041014754885176; 10 = 1.041815761046757; 11 = 1.055502933010359; 12 = 1.048454345559992; 13 = 1.060699658326322; 14 = 1.01770531771306 }
    14 = @{ 2 = 1.02; 3 = 1.035446138199214; 4 = 1.052120229389522; 5 = 1.045110329455921; 6 = 1.057394939979256; 9 = 1.041062977987075; 10 = 1.042018321345162; 11 = 1.055632239191205; 12 = 1.048647668521538; 13 = 1.060888134615797; 14 = 1.017774040384523 }
    15 = @{ 2 = 1.02; 3 = 1.035608905235249; 4 = 1.052219323094337; 5 = 1.045249200663493; 6 = 1.057530521096781; 9 = 1.041092615324374; 10 = 1.042143079800515; 11 = 1.055711822146015; 12 = 1.048766732773742; 13 = 1.061004189319167; 14 = 1.017816362495692 }
    16 = @{ 2 = 1.02; 3 = 1.036555700279217; 4 = 1.052794998312345; 5 = 1.046056953695081; 6 = 1.058318771113626; 9 = 1.041263980649551; 10 = 1.042868535923727; 11 = 1.056173705307466; 12 = 1.049459003298233; 13 = 1.061678585893958; 14 = 1.018062388682768 }
    17 = @{ 2 = 1.02; 3 = 1.037149096519543; 4 = 1.053155143904759; 5 = 1.046563167470239; 6 = 1.058812439501467; 9 = 1.04137047328609; 10 = 1.043322990984998; 11 = 1.056462268514982; 12 = 1.04989260360241; 13 = 1.062100659323513; 14 = 1.018216445890616 }
    18 = @{ 2 = 1.02; 3 = 1.037495035918751; 4 = 1.053364865600988; 5 = 1.046858266564651; 6 = 1.059100108559394; 9 = 1.041432228504937; 10 = 1.043587851118747; 11 = 1.056630163980171; 12 = 1.050145285488473; 13 = 1.062346504042754; 14 = 1.018306208747212 }
    19 = @{ 2 = 1.02; 3 = 1.037612962282902; 4 = 1.05343631682487; 5 = 1.046958859613972; 6 = 1.059198149052581; 9 = 1.0414532243868; 10 = 1.043678125160897; 11 = 1.056687341079351; 12 = 1.05023140466446; 13 = 1.062430272610163; 14 = 1.018336799322396 }
    20 = @{ 2 = 1.02; 3 = 1.03708544925379; 4 = 1.05311653941098; 5 = 1.046508872823642; 6 = 1.05875950249048; 9 = 1.041359084892561; 10 = 1.043274254625694; 11 = 1.05643135172844; 12 = 1.049846106206561; 13 = 1.062055410411778; 14 = 1.018199926965902 }
    21 = @{ 2 = 1.02; 3 = 1.035368334866553; 4 = 1.052072849135777; 5 = 1.045043947685697; 6 = 1.05733012444647; 9 = 1.041048793009209; 10 = 1.041958681903025; 11 = 1.055594179841925; 12 = 1.048590749821774; 13 = 1.060832648004248; 14 = 1.017753807481454 }
    22 = @{ 2 = 1.02; 3 = 1.034287048211007; 4 = 1.051413510258253; 5 = 1.044121346106947; 6 = 1.056428863390099; 9 = 1.040850447931863; 10 = 1.041129541860778; 11 = 1.055064026758866; 12 = 1.047799348474138; 13 = 1.060060719560959; 14 = 1.017472434375117 }
    23 = @{ 2 = 1.02; 3 = 1.034860419810664; 4 = 1.051763335887383; 5 = 1.044610583155584; 6 = 1.056906882888361; 9 = 1.04095590307249; 10 = 1.041569275149479; 11 = 1.055345431110362; 12 = 1.048219087079981; 13 = 1.060470231705221; 14 = 1.017621679698923 }
    24 = @{ 2 = 1.02; 3 = 1.037114209265608; 4 = 1.053133984186874; 5 = 1.046533406753196; 6 = 1.058783423310703; 9 = 1.041364231930374; 10 = 1.043296277155744; 11 = 1.056445322987196; 12 = 1.04986711708875; 13 = 1.062075857509245; 14 = 1.018207391453993 }
    25 = @{ 2 = 1.02; 3 = 1.039722630743139; 4 = 1.054711015752663; 5 = 1.048758212693281; 6 = 1.060950103685747; 9 = 1.041823951027067; 10 = 1.045291922279703; 11 = 1.057705275493183; 12 = 1.051770567219223; 13 = 1.063925635853634; 14 = 1.018883314388184 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Cells.Item($row, $col).Value = $data[$row][$col]
    }
}

Write-Output "Updated $($data.Count) rows of vm_pu results"